$wb = $excel.ActiveWorkbook

$ts = "2025-11-16 03:06:38"

$ws2 = $wb.Worksheets.Item("Главные")
$ws2.Range("C3").Value = 25
$ws2.Range("D3").Value = 439
$ws2.Range("E3").Value = 197
$ws2.Range("F3").Value = 242
$ws2.Range("G3").Value = 17.56
$ws2.Range("I3").Value = 9.68
$ws2.Range("J3").Value = 96
$ws2.Range("K3").Value = 101

$ws2.Range("C5").Value = 25
$ws2.Range("D5").Value = 411
$ws2.Range("E5").Value = 218
$ws2.Range("F5").Value = 193
$ws2.Range("G5").Value = 16.44
$ws2.Range("H5").Value = 8.72
$ws2.Range("I5").Value = 7.72
$ws2.Range("J5").Value = 104
$ws2.Range("K5").Value = 89

$ws2.Range("C10").Value = 17
$ws2.Range("D10").Value = 302
$ws2.Range("E10").Value = 156
$ws2.Range("F10").Value = 146
$ws2.Range("G10").Value = 17.76
$ws2.Range("H10").Value = 9.18
$ws2.Range("I10").Value = 8.59
$ws2.Range("J10").Value = 78
$ws2.Range("K10").Value = 63

$ws2.Range("C12").Value = 16
$ws2.Range("D12").Value = 281
$ws2.Range("E12").Value = 119
$ws2.Range("F12").Value = 162
$ws2.Range("G12").Value = 17.56
$ws2.Range("H12").Value = 7.44
$ws2.Range("I12").Value = 10.13
$ws2.Range("J12").Value = 47
$ws2.Range("K12").Value = 56
$ws2.Range("L12").Value = 3
$ws2.Range("V12").Value = 8

$ws2.Range("C13").Value = 10
$ws2.Range("D13").Value = 155
$ws2.Range("E13").Value = 88
$ws2.Range("F13").Value = 67
$ws2.Range("G13").Value = 15.5
$ws2.Range("H13").Value = 8.8
$ws2.Range("I13").Value = 6.7
$ws2.Range("J13").Value = 44
$ws2.Range("K13").Value = 31

$ws2.Range("C14").Value = 17
$ws2.Range("D14").Value = 211
$ws2.Range("E14").Value = 113
$ws2.Range("F14").Value = 98
$ws2.Range("G14").Value = 12.41
$ws2.Range("H14").Value = 6.65
$ws2.Range("I14").Value = 5.76
$ws2.Range("J14").Value = 54
$ws2.Range("K14").Value = 44
$ws2.Range("L14").Value = 1
$ws2.Range("V14").Value = 8

for ($r = 2; $r -le 26; $r++) {
    $ws2.Range("AA$r").Value = $ts
}

$ws3 = $wb.Worksheets.Item("Линейные")
$ws3.Range("C5").Value = 12
$ws3.Range("D5").Value = 170
$ws3.Range("E5").Value = 92
$ws3.Range("F5").Value = 78
$ws3.Range("G5").Value = 14.17
$ws3.Range("H5").Value = 7.67
$ws3.Range("I5").Value = 6.5
$ws3.Range("J5").Value = 46
$ws3.Range("K5").Value = 39

$ws3.Range("C15").Value = 21
$ws3.Range("D15").Value = 413
$ws3.Range("E15").Value = 219
$ws3.Range("F15").Value = 194
$ws3.Range("G15").Value = 19.67
$ws3.Range("H15").Value = 10.43
$ws3.Range("I15").Value = 9.24
$ws3.Range("J15").Value = 87
$ws3.Range("K15").Value = 77

$ws3.Range("C19").Value = 23
$ws3.Range("D19").Value = 403
$ws3.Range("E19").Value = 194
$ws3.Range("F19").Value = 209
$ws3.Range("G19").Value = 17.52
$ws3.Range("H19").Value = 8.43
$ws3.Range("I19").Value = 9.09
$ws3.Range("J19").Value = 92
$ws3.Range("K19").Value = 92

$ws3.Range("C21").Value = 27
$ws3.Range("D21").Value = 544
$ws3.Range("E21").Value = 230
$ws3.Range("F21").Value = 314
$ws3.Range("G21").Value = 20.15
$ws3.Range("H21").Value = 8.52
$ws3.Range("I21").Value = 11.63
$ws3.Range("J21").Value = 105
$ws3.Range("K21").Value = 127
$ws3.Range("L21").Value = 2
$ws3.Range("V21").Value = 12

$ws3.Range("C24").Value = 26
$ws3.Range("D24").Value = 475
$ws3.Range("E24").Value = 191
$ws3.Range("F24").Value = 284
$ws3.Range("G24").Value = 18.27
$ws3.Range("H24").Value = 7.35
$ws3.Range("I24").Value = 10.92
$ws3.Range("J24").Value = 83
$ws3.Range("K24").Value = 107
$ws3.Range("L24").Value = 3
$ws3.Range("V24").Value = 12

$ws3.Range("C26").Value = 22
$ws3.Range("D26").Value = 447
$ws3.Range("E26").Value = 191
$ws3.Range("F26").Value = 256
$ws3.Range("G26").Value = 20.32
$ws3.Range("H26").Value = 8.68
$ws3.Range("I26").Value = 11.64
$ws3.Range("J26").Value = 73
$ws3.Range("K26").Value = 73

for ($r = 2; $r -le 26; $r++) {
    $ws3.Range("AA$r").Value = $ts
}
